$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Parallel arrays: row numbers, new DAMSLTag (col I), new DialogAct (col J)
$rows = @(5, 11, 13, 16, 47, 55, 60, 67, 69, 70, 73, 75, 87, 88, 91, 114, 115, 117, 131, 135, 138, 144, 146, 150, 151, 153, 159, 164, 170, 185, 187, 195, 199, 204, 227, 232, 242, 244, 254, 264, 265, 266, 270, 283, 285, 291, 293, 294, 303, 306, 307, 311, 318, 322, 327, 335, 336, 337, 338, 339, 340, 343, 345, 352, 359, 364, 368, 373, 378, 392, 395, 402, 417, 420, 422)
$ivals = @("ba", "sd", "sd", "%", "b", "ba", "aa", "ba", "sv", "b", "b", "b", "b", "%", "%", "aa", "aa", "aa", "sd", "sd", "aa", "b", "sv", "sd", "sv", "sv", "sd", "aa", "aa", "sd", "sv", "aa", "b", "b", "aa", "%", "sd", "b", "b", "ba", "ba", "b", "sd", "sd", "%", "b", "aa", "sv", "sd", "ba", "b", "ba", "b", "aa", "aa", "sv", "ba", "sd", "sd", "b", "sd", "b", "aa", "sd", "ba", "aa", "sd", "aa", "ba", "b", "b", "ba", "aa", "sv", "b")
$jvals = @("Appreciation", "Statement-non-opinion", "Statement-non-opinion", "Uninterpretable", "Acknowledge (Backchannel)", "Appreciation", "Agree/Accept", "Appreciation", "Statement-opinion", "Acknowledge (Backchannel)", "Acknowledge (Backchannel)", "Acknowledge (Backchannel)", "Acknowledge (Backchannel)", "Uninterpretable", "Uninterpretable", "Agree/Accept", "Agree/Accept", "Agree/Accept", "Statement-non-opinion", "Statement-non-opinion", "Agree/Accept", "Acknowledge (Backchannel)", "Statement-opinion", "Statement-non-opinion", "Statement-opinion", "Statement-opinion", "Statement-non-opinion", "Agree/Accept", "Agree/Accept", "Statement-non-opinion", "Statement-opinion", "Agree/Accept", "Acknowledge (Backchannel)", "Acknowledge (Backchannel)", "Agree/Accept", "Uninterpretable", "Statement-non-opinion", "Acknowledge (Backchannel)", "Acknowledge (Backchannel)", "Appreciation", "Appreciation", "Acknowledge (Backchannel)", "Statement-non-opinion", "Statement-non-opinion", "Uninterpretable", "Acknowledge (Backchannel)", "Agree/Accept", "Statement-opinion", "Statement-non-opinion", "Appreciation", "Acknowledge (Backchannel)", "Appreciation", "Acknowledge (Backchannel)", "Agree/Accept", "Agree/Accept", "Statement-opinion", "Appreciation", "Statement-non-opinion", "Statement-non-opinion", "Acknowledge (Backchannel)", "Statement-non-opinion", "Acknowledge (Backchannel)", "Agree/Accept", "Statement-non-opinion", "Appreciation", "Agree/Accept", "Statement-non-opinion", "Agree/Accept", "Appreciation", "Acknowledge (Backchannel)", "Acknowledge (Backchannel)", "Appreciation", "Agree/Accept", "Statement-opinion", "Acknowledge (Backchannel)")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 9).Value = $ivals[$i]
    $ws.Cells.Item($r, 10).Value = $jvals[$i]
}
